# Applies a cyclic rotation of the data rows 2-4 on Sheet1:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# Only the columns that actually differ between the rotated rows need to be
# written (the other columns already hold identical values across the rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture the "before" values for the columns that change under the rotation.
# NOTE: use Value() (method call syntax) rather than bare .Value when *reading*
# - in this runtime the bare property getter does not evaluate correctly.
$row2 = @{
    D = $ws.Range("D2").Value()
    M = $ws.Range("M2").Value()
    N = $ws.Range("N2").Value()
    O = $ws.Range("O2").Value()
    P = $ws.Range("P2").Value()
    R = $ws.Range("R2").Value()
    S = $ws.Range("S2").Value()
}

$row3 = @{
    D = $ws.Range("D3").Value()
    M = $ws.Range("M3").Value()
    N = $ws.Range("N3").Value()
    O = $ws.Range("O3").Value()
    P = $ws.Range("P3").Value()
    R = $ws.Range("R3").Value()
    S = $ws.Range("S3").Value()
}

$row4 = @{
    D = $ws.Range("D4").Value()
    M = $ws.Range("M4").Value()
    N = $ws.Range("N4").Value()
    O = $ws.Range("O4").Value()
    P = $ws.Range("P4").Value()
    R = $ws.Range("R4").Value()
    S = $ws.Range("S4").Value()
}

# Write row 2 <- old row 4
$ws.Range("D2").Value = $row4.D
$ws.Range("M2").Value = $row4.M
$ws.Range("N2").Value = $row4.N
$ws.Range("O2").Value = $row4.O
$ws.Range("P2").Value = $row4.P
$ws.Range("R2").Value = $row4.R
$ws.Range("S2").Value = $row4.S

# Write row 3 <- old row 2
$ws.Range("D3").Value = $row2.D
$ws.Range("M3").Value = $row2.M
$ws.Range("N3").Value = $row2.N
$ws.Range("O3").Value = $row2.O
$ws.Range("P3").Value = $row2.P
$ws.Range("R3").Value = $row2.R
$ws.Range("S3").Value = $row2.S

# Write row 4 <- old row 3
$ws.Range("D4").Value = $row3.D
$ws.Range("M4").Value = $row3.M
$ws.Range("N4").Value = $row3.N
$ws.Range("O4").Value = $row3.O
$ws.Range("P4").Value = $row3.P
$ws.Range("R4").Value = $row3.R
$ws.Range("S4").Value = $row3.S
